$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 545.9
$ws.Range("I2").Value = 484.33334
$ws.Range("J2").Value = 1100
$ws.Range("K2").Value = 484.33334
$ws.Range("L2").Value = 1100
$ws.Range("M2").Value = -371.33334
$ws.Range("N2").Value = -1326
$ws.Range("H17").Value = 853.9506
$ws.Range("I17").Value = 1444.4445
$ws.Range("J17").Value = 780.1389
$ws.Range("K17").Value = 4333.333500000001
$ws.Range("L17").Value = 2340.4167
$ws.Range("M17").Value = -4165.333500000001
$ws.Range("N17").Value = -2676.4167
$ws.Range("H92").Value = 1351.826
$ws.Range("I92").Value = 1215.1052
$ws.Range("J92").Value = 2001.25
$ws.Range("K92").Value = 1215.1052
$ws.Range("L92").Value = 2001.25
$ws.Range("M92").Value = 32.89480000000003
$ws.Range("N92").Value = -4497.25
$ws.Range("H94").Value = 1974.75
$ws.Range("I94").Value = 1974.75
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1974.75
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1523.75
$ws.Range("N94").ClearContents()
$ws.Range("H100").Value = 18185000
$ws.Range("I100").Value = 20003160
$ws.Range("J100").Value = 3400
$ws.Range("K100").Value = 20003160
$ws.Range("L100").Value = 3400
$ws.Range("M100").Value = -20002619
$ws.Range("N100").Value = -4482
$ws.Range("H103").Value = 9096.4
$ws.Range("I103").Value = 670.625
$ws.Range("J103").Value = 42799.5
$ws.Range("K103").Value = 2011.875
$ws.Range("L103").Value = 128398.5
$ws.Range("M103").Value = -1425.875
$ws.Range("N103").Value = -129570.5
$ws.Range("H106").Value = 2142.5
$ws.Range("I106").Value = 1332.9166
$ws.Range("J106").Value = 7000
$ws.Range("K106").Value = 1332.9166
$ws.Range("L106").Value = 7000
$ws.Range("M106").Value = -701.9166
$ws.Range("N106").Value = -8262
$ws.Range("H133").Value = 46157.273
$ws.Range("J133").Value = 46157.273
$ws.Range("L133").Value = 46157.273
$ws.Range("N133").Value = -56277.273
$ws.Range("H138").Value = 3299.93
$ws.Range("I138").Value = 727.9429
$ws.Range("J138").Value = 4684.846
$ws.Range("K138").Value = 2183.8287
$ws.Range("L138").Value = 14054.538
$ws.Range("M138").Value = 2956.1713
$ws.Range("N138").Value = -24334.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 38394
$ws.Range("J76").Value = 38394
$ws.Range("L76").Value = 38394
$ws.Range("N76").Value = -39070
$ws.Range("H79").Value = 38394
$ws.Range("J79").Value = 38394
$ws.Range("L79").Value = 38394
$ws.Range("N79").Value = -40734
$ws.Range("H97").Value = 927.8570999999999
$ws.Range("I97").Value = 914.25
$ws.Range("K97").Value = 914.25
$ws.Range("M97").Value = -418.25
$ws.Range("H102").Value = 1253.4783
$ws.Range("I102").Value = 1084.1177
$ws.Range("J102").Value = 1733.3334
$ws.Range("K102").Value = 1084.1177
$ws.Range("L102").Value = 1733.3334
$ws.Range("M102").Value = 537.8823
$ws.Range("N102").Value = -4977.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1156.8
$ws.Range("I94").Value = 1271
$ws.Range("J94").Value = 700
$ws.Range("K94").Value = 1271
$ws.Range("L94").Value = 700
$ws.Range("M94").Value = -820
$ws.Range("N94").Value = -1602
$ws.Range("H99").Value = 1698.1364
$ws.Range("I99").Value = 1149.1333
$ws.Range("J99").Value = 2874.5715
$ws.Range("K99").Value = 1149.1333
$ws.Range("L99").Value = 2874.5715
$ws.Range("M99").Value = 348.8667
$ws.Range("N99").Value = -5870.5715
$ws.Range("H105").Value = 1695.7307
$ws.Range("I105").Value = 1656.826
$ws.Range("K105").Value = 1656.826
$ws.Range("M105").Value = 90.17399999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 56260.145
$ws.Range("J68").Value = 56260.145
$ws.Range("L68").Value = 56260.145
$ws.Range("N68").Value = -57758.145
$ws.Range("H71").Value = 56260.145
$ws.Range("J71").Value = 56260.145
$ws.Range("L71").Value = 168780.435
$ws.Range("N71").Value = -176268.435
$ws.Range("H106").Value = 34699.5
$ws.Range("J106").Value = 34699.5
$ws.Range("L106").Value = 34699.5
$ws.Range("N106").Value = -37223.5
$ws.Range("H107").Value = 826.0357
$ws.Range("I107").Value = 552.6667
$ws.Range("J107").Value = 1141.4615
$ws.Range("K107").Value = 552.6667
$ws.Range("L107").Value = 1141.4615
$ws.Range("M107").Value = 1367.3333
$ws.Range("N107").Value = -4981.461499999999
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H123").Value = 36602.375
$ws.Range("J123").Value = 36602.375
$ws.Range("L123").Value = 36602.375
$ws.Range("N123").Value = -46402.375
$ws.Range("H134").Value = 1350
$ws.Range("I134").Value = 590.625
$ws.Range("J134").Value = 3780
$ws.Range("K134").Value = 1771.875
$ws.Range("L134").Value = 11340
$ws.Range("M134").Value = 763.125
$ws.Range("N134").Value = -16410

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 2649.8572
$ws.Range("I21").Value = 633.3333
$ws.Range("J21").Value = 3199.818
$ws.Range("K21").Value = 1899.9999
$ws.Range("L21").Value = 9599.454000000002
$ws.Range("M21").Value = -1726.9999
$ws.Range("N21").Value = -9945.454000000002
$ws.Range("H113").Value = 521.8946999999999
$ws.Range("I113").Value = 496.7
$ws.Range("J113").Value = 549.8889
$ws.Range("K113").Value = 1490.1
$ws.Range("L113").Value = 1649.6667
$ws.Range("M113").Value = 679.9000000000001
$ws.Range("N113").Value = -5989.6667
$ws.Range("H121").Value = 2627.3333
$ws.Range("I121").Value = 210
$ws.Range("J121").Value = 2754.5615
$ws.Range("K121").Value = 630
$ws.Range("L121").Value = 8263.684499999999
$ws.Range("M121").Value = 680
$ws.Range("N121").Value = -10883.6845
$ws.Range("H131").Value = 6579747
$ws.Range("I131").Value = 62500310
$ws.Range("J131").Value = 857.0294
$ws.Range("K131").Value = 187500930
$ws.Range("L131").Value = 2571.0882
$ws.Range("M131").Value = -187495890
$ws.Range("N131").Value = -12651.0882

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 15929.8
$ws.Range("J39").Value = 15929.8
$ws.Range("L39").Value = 15929.8
$ws.Range("N39").Value = -16993.8
$ws.Range("H113").Value = 1966.4445
$ws.Range("I113").Value = 1899.7142
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 1899.7142
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = 270.2858000000001
$ws.Range("N113").Value = -6540
$ws.Range("H132").Value = 1835.4637
$ws.Range("I132").Value = 1202.6666
$ws.Range("K132").Value = 3607.9998
$ws.Range("M132").Value = -1077.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 5877
$ws.Range("I23").Value = 506
$ws.Range("J23").Value = 11248
$ws.Range("K23").Value = 506
$ws.Range("L23").Value = 11248
$ws.Range("M23").Value = -276
$ws.Range("N23").Value = -11708
$ws.Range("H54").Value = 35026.332
$ws.Range("J54").Value = 35026.332
$ws.Range("L54").Value = 35026.332
$ws.Range("N54").Value = -36314.332
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360
$ws.Range("H82").Value = 1235.3903
$ws.Range("I82").Value = 602.4091
$ws.Range("J82").Value = 1968.3158
$ws.Range("K82").Value = 602.4091
$ws.Range("L82").Value = 1968.3158
$ws.Range("M82").Value = -241.4091
$ws.Range("N82").Value = -2690.3158
$ws.Range("H85").Value = 1235.3903
$ws.Range("I85").Value = 602.4091
$ws.Range("J85").Value = 1968.3158
$ws.Range("K85").Value = 602.4091
$ws.Range("L85").Value = 1968.3158
$ws.Range("M85").Value = 645.5909
$ws.Range("N85").Value = -4464.3158
$ws.Range("H93").Value = 4274829
$ws.Range("I93").Value = 6536796.5
$ws.Range("J93").Value = 2224.2222
$ws.Range("K93").Value = 6536796.5
$ws.Range("L93").Value = 2224.2222
$ws.Range("M93").Value = -6535548.5
$ws.Range("N93").Value = -4720.2222
$ws.Range("H100").Value = 1733.3334
$ws.Range("I100").Value = 1371.4286
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1371.4286
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -830.4286
$ws.Range("N100").Value = -4082
$ws.Range("H122").Value = 4270.5293
$ws.Range("I122").Value = 1914.2858
$ws.Range("J122").Value = 5919.9
$ws.Range("K122").Value = 5742.857400000001
$ws.Range("L122").Value = 17759.7
$ws.Range("M122").Value = -3292.857400000001
$ws.Range("N122").Value = -22659.7
$ws.Range("H123").Value = 27651.25
$ws.Range("J123").Value = 27651.25
$ws.Range("L123").Value = 27651.25
$ws.Range("N123").Value = -37451.25
$ws.Range("H136").Value = 2204.2424
$ws.Range("I136").Value = 1176.4286
$ws.Range("K136").Value = 3529.2858
$ws.Range("M136").Value = -979.2857999999997

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 24260
$ws.Range("J57").Value = 24260
$ws.Range("L57").Value = 24260
$ws.Range("N57").Value = -25768
$ws.Range("H96").Value = 94038450
$ws.Range("J96").Value = 3969651.2
$ws.Range("L96").Value = 3969651.2
$ws.Range("N96").Value = -3972397.2
$ws.Range("H100").Value = 811.7692
$ws.Range("I100").Value = 671.1667
$ws.Range("J100").Value = 2499
$ws.Range("K100").Value = 1342.3334
$ws.Range("L100").Value = 4998
$ws.Range("M100").Value = -801.3334
$ws.Range("N100").Value = -6080
